$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped price/volume/coin updates cell-by-cell.
# Numeric-looking text (e.g. "591.97") is forced to stay text: Excel
# auto-converts a plain numeric literal assigned via .Value into a real
# number, so we flip the cell to text format first, assign, then restore
# the original (default) cell style so no stray formatting is introduced.

$ws.Cells.Item(2, 4).Value = '66.622.05'
$ws.Cells.Item(2, 5).Value = '  +1.33%  '

$ws.Cells.Item(3, 4).Value = '3.495.08'
$ws.Cells.Item(3, 5).Value = '  +0.78%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '591.97'
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  +1.79%  '

$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '168.67'
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  +0.46%  '

$ws.Cells.Item(7, 5).Value = '  +0.04%  '

$ws.Cells.Item(8, 2).Value = 'XRP'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$origStyle = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.592'
$ws.Cells.Item(8, 4).Style = $origStyle
$ws.Cells.Item(8, 5).Value = '  +4.93%  '

$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.128'
$ws.Cells.Item(9, 4).Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  +4.45%  '

$ws.Cells.Item(10, 2).Value = 'Toncoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.32'
$ws.Cells.Item(10, 4).Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  +0.16%  '

$ws.Cells.Item(11, 2).Value = 'Cardano'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$origStyle = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.431'
$ws.Cells.Item(11, 4).Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  -0.04%  '

$ws.Cells.Item(12, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(12, 4).Value = '4.098.20'
$ws.Cells.Item(12, 5).Value = '  +1.15%  '

$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$origStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.134'
$ws.Cells.Item(13, 4).Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  -0.31%  '

$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '28.14'
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  +2.22%  '

$ws.Cells.Item(15, 2).Value = 'ShibaInu'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$origStyle = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000178'
$ws.Cells.Item(15, 4).Style = $origStyle
$ws.Cells.Item(15, 5).Value = '  +0.98%  '

$ws.Cells.Item(16, 4).Value = '66.634.49'
$ws.Cells.Item(16, 5).Value = '  +1.45%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '3.486.60'
$ws.Cells.Item(17, 5).Value = '  +1.31%  '

$ws.Cells.Item(18, 2).Value = 'Polkadot'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$origStyle = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.29'
$ws.Cells.Item(18, 4).Style = $origStyle
$ws.Cells.Item(18, 5).Value = '  +0.80%  '

$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '14.04'
$ws.Cells.Item(19, 4).Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  +1.89%  '

$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '393.26'
$ws.Cells.Item(20, 4).Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  +2.34%  '

$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$origStyle = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '7.91'
$ws.Cells.Item(21, 4).Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  -0.60%  '

$ws.Cells.Item(22, 2).Value = 'Litecoin'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$origStyle = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '73.00'
$ws.Cells.Item(22, 4).Style = $origStyle
$ws.Cells.Item(22, 5).Value = '  +1.73%  '

$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  -0.13%  '

$ws.Cells.Item(24, 2).Value = 'Polygon'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.534'
$ws.Cells.Item(24, 4).Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  +2.52%  '

$ws.Cells.Item(25, 2).Value = 'PEPE'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.0000122'
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  +0.91%  '

$ws.Cells.Item(26, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '10.18'
$ws.Cells.Item(26, 4).Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  +3.30%  '

$ws.Cells.Item(27, 2).Value = 'Kaspa'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.180'
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  -0.96%  '

$ws.Cells.Item(28, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$origStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.999'
$ws.Cells.Item(28, 4).Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  +0.04%  '

$ws.Cells.Item(29, 2).Value = 'NEARProtocol'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$origStyle = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '6.37'
$ws.Cells.Item(29, 4).Style = $origStyle
$ws.Cells.Item(29, 5).Value = '  +2.04%  '

$ws.Cells.Item(30, 2).Value = 'Fetch.AI'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.45'
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  -0.14%  '

$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.05'
$ws.Cells.Item(31, 4).Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  +1.49%  '

$ws.Cells.Item(32, 2).Value = 'EthereumClassic'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '23.68'
$ws.Cells.Item(32, 4).Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  +1.54%  '

$ws.Cells.Item(33, 2).Value = 'Aptos'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$origStyle = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '7.33'
$ws.Cells.Item(33, 4).Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  +0.26%  '

$ws.Cells.Item(34, 2).Value = 'USDe'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.00'
$ws.Cells.Item(34, 4).Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  +0.03%  '

$origStyle = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.59'
$ws.Cells.Item(35, 4).Style = $origStyle
$ws.Cells.Item(35, 5).Value = '  +4.50%  '

$origStyle = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '162.41'
$ws.Cells.Item(36, 4).Style = $origStyle
$ws.Cells.Item(36, 5).Value = '  +1.39%  '

$origStyle = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.895'
$ws.Cells.Item(37, 4).Style = $origStyle
$ws.Cells.Item(37, 5).Value = '  +0.37%  '

$origStyle = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.92'
$ws.Cells.Item(38, 4).Style = $origStyle
$ws.Cells.Item(38, 5).Value = '  +1.84%  '

$ws.Cells.Item(39, 5).Value = '  +2.24%  '

$origStyle = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '4.64'
$ws.Cells.Item(40, 4).Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  +3.86%  '

$origStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0739'
$ws.Cells.Item(41, 4).Style = $origStyle
$ws.Cells.Item(41, 5).Value = '  +0.37%  '

$origStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '26.44'
$ws.Cells.Item(42, 4).Style = $origStyle
$ws.Cells.Item(42, 5).Value = '  +0.69%  '

$origStyle = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '26.80'
$ws.Cells.Item(43, 4).Style = $origStyle
$ws.Cells.Item(43, 5).Value = '  +0.17%  '

$ws.Cells.Item(44, 4).Value = '2.771.69'
$ws.Cells.Item(44, 5).Value = '  -1.40%  '

$origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '42.89'
$ws.Cells.Item(45, 4).Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  -0.60%  '

$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.55'
$ws.Cells.Item(46, 4).Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  +2.24%  '

$origStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0310'
$ws.Cells.Item(47, 4).Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  -0.12%  '

$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '343.89'
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  +1.62%  '

$origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.09'
$ws.Cells.Item(49, 4).Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  +1.35%  '

$origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '33.97'
$ws.Cells.Item(50, 4).Style = $origStyle
$ws.Cells.Item(50, 5).Value = '  +4.65%  '

$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.856'
$ws.Cells.Item(51, 4).Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  +2.73%  '
